$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" / "Volume(1h)" columns with the latest scrape, and
# (rows 20/21) swap Polkadot and BitcoinCash back into rank order.
#
# A few Price values (D13, D40, D45 below) end in a significant trailing
# zero (e.g. "0.340"). A bare Range.Value assignment lets Excel's normal
# text/number autodetection kick in and it would silently collapse that to
# the number 0.34, dropping the zero - so those three are entered with a
# leading apostrophe, same as typing them into Excel by hand, to force
# literal text and keep the exact digits from the source data.

$ws.Range("D2").Value = "58.056.81"
$ws.Range("E2").Value = "  -2.00%  "

$ws.Range("D3").Value = "2.469.47"
$ws.Range("E3").Value = "  -2.25%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "517.91"
$ws.Range("E5").Value = "  -3.52%  "

$ws.Range("D6").Value = "131.61"
$ws.Range("E6").Value = "  -3.94%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "0.558"
$ws.Range("E8").Value = "  -1.90%  "

$ws.Range("D9").Value = "2.472.20"
$ws.Range("E9").Value = "  -2.08%  "

$ws.Range("D10").Value = "0.0992"
$ws.Range("E10").Value = "  -2.26%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("D12").Value = "5.34"
$ws.Range("E12").Value = "  +0.80%  "

$ws.Range("D13").Value = "'0.340"
$ws.Range("E13").Value = "  -2.38%  "

$ws.Range("D14").Value = "2.907.57"
$ws.Range("E14").Value = "  -2.24%  "

$ws.Range("D15").Value = "58.007.32"
$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("D16").Value = "22.38"
$ws.Range("E16").Value = "  -3.12%  "

$ws.Range("E17").Value = "  -1.98%  "

$ws.Range("D18").Value = "2.469.31"
$ws.Range("E18").Value = "  -2.38%  "

$ws.Range("D19").Value = "10.74"
$ws.Range("E19").Value = "  -3.77%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "4.18"
$ws.Range("E20").Value = "  -2.49%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "320.41"
$ws.Range("E21").Value = "  -1.14%  "

$ws.Range("D23").Value = "5.72"
$ws.Range("E23").Value = "  -4.22%  "

$ws.Range("D24").Value = "64.17"
$ws.Range("E24").Value = "  -2.01%  "

$ws.Range("E25").Value = "  -2.62%  "

$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("E27").Value = "  -3.47%  "

$ws.Range("D28").Value = "7.33"
$ws.Range("E28").Value = "  -2.81%  "

$ws.Range("D29").Value = "0.0₃0750"
$ws.Range("E29").Value = "  -3.04%  "

$ws.Range("E30").Value = "  -4.45%  "

$ws.Range("D31").Value = "165.53"
$ws.Range("E31").Value = "  -2.88%  "

$ws.Range("D32").Value = "6.25"
$ws.Range("E32").Value = "  -6.65%  "

$ws.Range("D33").Value = "1.17"
$ws.Range("E33").Value = "  -1.80%  "

$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").Value = "18.03"
$ws.Range("E36").Value = "  -1.96%  "

$ws.Range("E37").Value = "  -8.01%  "

$ws.Range("D38").Value = "3.98"
$ws.Range("E38").Value = "  -2.95%  "

$ws.Range("E39").Value = "  -4.78%  "

$ws.Range("D40").Value = "'0.790"
$ws.Range("E40").Value = "  -2.85%  "

$ws.Range("E41").Value = "  -4.44%  "

$ws.Range("D42").Value = "272.16"
$ws.Range("E42").Value = "  -5.08%  "

$ws.Range("D43").Value = "4.98"
$ws.Range("E43").Value = "  -2.68%  "

$ws.Range("D44").Value = "0.591"
$ws.Range("E44").Value = "  -3.01%  "

$ws.Range("D45").Value = "'126.80"
$ws.Range("E45").Value = "  -3.41%  "

$ws.Range("D46").Value = "0.0905"
$ws.Range("E46").Value = "  -2.04%  "

$ws.Range("D47").Value = "0.0487"
$ws.Range("E47").Value = "  -4.07%  "

$ws.Range("E48").Value = "  -4.23%  "

$ws.Range("D49").Value = "16.94"
$ws.Range("E49").Value = "  -2.79%  "

$ws.Range("D50").Value = "1.725.57"
$ws.Range("E50").Value = "  -2.11%  "

$ws.Range("D51").Value = "0.973"
$ws.Range("E51").Value = "  -1.77%  "
